# Regenerate the handback report: the handoff file pair changed from
#   416ce689-5ba7-4604-92bd-dd0b924fa3fe.md / 4ffe8906-cf1d-4fba-866d-8e285fb1ae43.md
# to
#   6ed8456b-c2ad-4330-9572-caf6d7fa80c0.md / fffff5089c33-728a-49b8-bee3-dc59b445d663.md
# with new handoff/handback xlf file names and new timestamps.

$wb = $excel.ActiveWorkbook

$newA = "6ed8456b-c2ad-4330-9572-caf6d7fa80c0.md"
$newB = "fffff5089c33-728a-49b8-bee3-dc59b445d663.md"

$newXlfZh  = "6ed8456b-c2ad-4330-9572-caf6d7fa80c0.41d923ced29b8f299034b77d0df713481a4d1485.zh-cn.xlf"
$newXlfDe  = "6ed8456b-c2ad-4330-9572-caf6d7fa80c0.41d923ced29b8f299034b77d0df713481a4d1485.de-de.xlf"

$newZhHandoffTime  = "2016-03-18 07:29:03"
$newZhHandbackTime = "2016-03-18 07:29:19"
$newDeHandoffTime  = "2016-03-18 07:29:06"
$newDeHandbackTime = "2016-03-18 07:29:24"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newA
$wsOverview.Range("A3").Value = $newB

foreach ($hl in $wsOverview.Hyperlinks) {
    $addr = $hl.Range.Address()
    switch ($addr) {
        '$A$2' { $hl.TextToDisplay = $newA }
        '$A$3' { $hl.TextToDisplay = $newB }
    }
}

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newA
$wsZh.Range("D2").Value = $newXlfZh
$wsZh.Range("E2").Value = $newZhHandoffTime
$wsZh.Range("F2").Value = $newA
$wsZh.Range("G2").Value = $newXlfZh
$wsZh.Range("H2").Value = $newZhHandbackTime

$wsZh.Range("A3").Value = $newB
$wsZh.Range("D3").Value = $newXlfZh
$wsZh.Range("E3").Value = $newZhHandoffTime
$wsZh.Range("F3").Value = $newB
$wsZh.Range("G3").Value = $newXlfZh
$wsZh.Range("H3").Value = $newZhHandbackTime

foreach ($hl in $wsZh.Hyperlinks) {
    $addr = $hl.Range.Address()
    switch ($addr) {
        '$A$2' { $hl.TextToDisplay = $newA }
        '$D$2' { $hl.TextToDisplay = $newXlfZh }
        '$F$2' { $hl.TextToDisplay = $newA }
        '$G$2' { $hl.TextToDisplay = $newXlfZh }
        '$A$3' { $hl.TextToDisplay = $newB }
        '$D$3' { $hl.TextToDisplay = $newXlfZh }
        '$F$3' { $hl.TextToDisplay = $newB }
        '$G$3' { $hl.TextToDisplay = $newXlfZh }
    }
}

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newA
$wsDe.Range("D2").Value = $newXlfDe
$wsDe.Range("E2").Value = $newDeHandoffTime
$wsDe.Range("F2").Value = $newA
$wsDe.Range("G2").Value = $newXlfDe
$wsDe.Range("H2").Value = $newDeHandbackTime

$wsDe.Range("A3").Value = $newB
$wsDe.Range("D3").Value = $newXlfDe
$wsDe.Range("E3").Value = $newDeHandoffTime
$wsDe.Range("F3").Value = $newB
$wsDe.Range("G3").Value = $newXlfDe
$wsDe.Range("H3").Value = $newDeHandbackTime

foreach ($hl in $wsDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    switch ($addr) {
        '$A$2' { $hl.TextToDisplay = $newA }
        '$D$2' { $hl.TextToDisplay = $newXlfDe }
        '$F$2' { $hl.TextToDisplay = $newA }
        '$G$2' { $hl.TextToDisplay = $newXlfDe }
        '$A$3' { $hl.TextToDisplay = $newB }
        '$D$3' { $hl.TextToDisplay = $newXlfDe }
        '$F$3' { $hl.TextToDisplay = $newB }
        '$G$3' { $hl.TextToDisplay = $newXlfDe }
    }
}
